$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.266.99"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3
$ws.Range("D3").Value = "3.202.88"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "607.33"
$ws.Range("E5").Value = "  +1.63%  "

# Row 6
$ws.Range("D6").Value = "156.00"
$ws.Range("E6").Value = "  +0.61%  "

# Row 8
$ws.Range("D8").Value = "3.202.72"
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  -1.67%  "

# Row 10
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.30%  "

# Row 11
$ws.Range("E11").Value = "  -4.36%  "

# Row 12
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  -3.28%  "

# Row 13
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  -0.64%  "

# Row 14
$ws.Range("D14").Value = "38.34"
$ws.Range("E14").Value = "  -2.48%  "

# Row 15
$ws.Range("D15").Value = "3.728.55"
$ws.Range("E15").Value = "  +0.34%  "

# Row 16
$ws.Range("D16").Value = "66.411.79"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").Value = "7.30"
$ws.Range("E17").Value = "  -2.89%  "

# Row 18
$ws.Range("D18").Value = "3.205.14"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("E19").Value = "  +1.29%  "

# Row 20
$ws.Range("D20").Value = "506.14"
$ws.Range("E20").Value = "  -2.69%  "

# Row 21
$ws.Range("D21").Value = "15.30"
$ws.Range("E21").Value = "  -0.99%  "

# Row 22
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  -1.46%  "

# Row 23
$ws.Range("E23").Value = "  -1.78%  "

# Row 24
$ws.Range("D24").Value = "14.59"
$ws.Range("E24").Value = "  -2.56%  "

# Row 25
$ws.Range("D25").Value = "85.10"
$ws.Range("E25").Value = "  -1.06%  "

# Row 26
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -0.48%  "

# Row 28
$ws.Range("D28").Value = "9.03"
$ws.Range("E28").Value = "  -2.72%  "

# Row 29
$ws.Range("E29").Value = "  -0.73%  "

# Row 30
$ws.Range("E30").Value = "  +41.59%  "

# Row 31
$ws.Range("E31").Value = "  -0.65%  "

# Row 32
$ws.Range("E32").Value = "  -1.86%  "

# Row 33
$ws.Range("D33").Value = "28.21"
$ws.Range("E33").Value = "  -0.59%  "

# Row 34
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("E35").Value = "  -5.32%  "

# Row 36
$ws.Range("D36").Value = "6.43"
$ws.Range("E36").Value = "  -1.91%  "

# Row 37
$ws.Range("D37").Value = "500.53"
$ws.Range("E37").Value = "  -2.01%  "

# Row 38
$ws.Range("D38").Value = "55.37"
$ws.Range("E38").Value = "  +0.88%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0767"
$ws.Range("E39").Value = "  +12.84%  "

# Row 40
$ws.Range("E40").Value = "  +2.04%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.05"
$ws.Range("E41").Value = "  +5.51%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0418"
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("E43").Value = "  -2.51%  "

# Row 44
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("D45").Value = "2.911.42"
$ws.Range("E45").Value = "  +0.34%  "

# Row 46
$ws.Range("E46").Value = "  -1.31%  "

# Row 47
$ws.Range("D47").Value = "28.14"
$ws.Range("E47").Value = "  -1.38%  "

# Row 48
$ws.Range("E48").Value = "  +1.82%  "

# Row 50
$ws.Range("E50").Value = "  -1.04%  "

# Row 51
$ws.Range("D51").Value = "122.08"
$ws.Range("E51").Value = "  +0.49%  "
